$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old layout so no stale cells remain ---
$ws.Range("A1:F20").ClearContents()

# --- Row 1: headers ---
$ws.Range("A1").Value = "Part"
$ws.Range("B1").Value = "Number"
$ws.Range("C1").Value = "Needed"
$ws.Range("D1").Value = "Order"
$ws.Range("E1").Value = "Unit Cost"

# --- Row 2: ATMEGA ---
$ws.Range("A2").Value = "ATMEGA"
$ws.Range("B2").Value = "ATMEGA32U4-AU"
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 4.12

# --- Row 3: RFM69HCW ---
$ws.Range("A3").Value = "RFM69HCW - 915MHz"
$ws.Range("B3").Value = "1568-1394-ND"
$ws.Range("C3").Value = 1

# --- Row 4: 3.3V regulator (was row 2) ---
$ws.Range("A4").Value = "3.3V regulator"
$ws.Range("B4").Value = "296-39452-1-ND"
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 1.23

# --- Row 5: 10uH inductor (was row 3) ---
$ws.Range("A5").Value = "10uH inductor"
$ws.Range("B5").Value = "587-2886-1-ND"
$ws.Range("C5").Value = 1
$ws.Range("E5").Value = 0.29

# --- Row 6: 47uF ceramic (was row 4, now only Part col) ---
$ws.Range("A6").Value = "47uF ceramic"

# --- Row 7: 4.7uF ceramic (was row 5) ---
$ws.Range("A7").Value = "4.7uF ceramic"
$ws.Range("B7").Value = "587-1780-1-ND"
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0.7

# --- Row 8: 1k resistor ---
$ws.Range("A8").Value = "1k resistor"
$ws.Range("B8").Value = "541-3991-1-ND"

# --- Row 9: 10k resistor ---
$ws.Range("A9").Value = "10k resistor"
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = "Have"

# --- Row 10: reset button ---
$ws.Range("A10").Value = "reset button"
$ws.Range("B10").Value = "401-1426-1-ND"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 0.52

# --- Row 11: 0.1uF ceramic ---
$ws.Range("A11").Value = "0.1uF ceramic"

# --- Row 12: 8MHz crystal ---
$ws.Range("A12").Value = "8MHz crystal"
$ws.Range("B12").Value = "535-10212-1-ND"
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 0.27

# --- Row 13: 18pF ceramic cap ---
$ws.Range("A13").Value = "18pF ceramic cap"
$ws.Range("C13").Value = 2

# --- Styles: small font (Arial 7pt black) on B7:B8 ---
$fontRange = $ws.Range("B7:B8")
$fontRange.Font.Name = "Arial"
$fontRange.Font.Size = 7
$fontRange.Font.Color = 0

# --- Styles: currency number format on the Unit Cost values ---
$currencyFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"
$ws.Range("E4").NumberFormat = $currencyFormat
$ws.Range("E5").NumberFormat = $currencyFormat
$ws.Range("E7").NumberFormat = $currencyFormat
$ws.Range("E10").NumberFormat = $currencyFormat
$ws.Range("E12").NumberFormat = $currencyFormat

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 20.666666666666668
$ws.Columns("E").ColumnWidth = 8

# --- Dangling external reference artifact (copy/paste from another workbook) ---
$ws.Range("Z1").Formula = "='[Book2.xlsx]Sheet1'!A1"
$ws.Range("Z1").Value = "x"
$ws.Range("Z1").ClearContents()

# --- Selection ---
$ws.Range("C13").Select()
